$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.69"
$ws.Range("E2").Value = "'2.27%"
$ws.Range("G2").Value = "'21"

$ws.Range("D3").Value = "'37.69"
$ws.Range("E3").Value = "'0.72%"
$ws.Range("G3").Value = "'21"

$ws.Range("D4").Value = "'5.166"
$ws.Range("E4").Value = "'1.46%"
$ws.Range("G4").Value = "'21"

$ws.Range("D5").Value = "'0.07886"
$ws.Range("E5").Value = "'2.07%"
$ws.Range("G5").Value = "'21"

$ws.Range("D6").Value = "'4.418"
$ws.Range("E6").Value = "'1.56%"
$ws.Range("G6").Value = "'21"

$ws.Range("E7").Value = "'1.94%"
$ws.Range("G7").Value = "'21"

$ws.Range("D8").Value = "'8.302"
$ws.Range("E8").Value = "'1.28%"
$ws.Range("G8").Value = "'21"

$ws.Range("D9").Value = "'2.853"
$ws.Range("G9").Value = "'21"

$ws.Range("D10").Value = "'0.9213"
$ws.Range("E10").Value = "'0.60%"
$ws.Range("G10").Value = "'21"

$ws.Range("E11").Value = "'5.75%"
$ws.Range("G11").Value = "'21"

$ws.Range("D12").Value = "'0.1931"
$ws.Range("E12").Value = "'2.75%"
$ws.Range("G12").Value = "'21"

$ws.Range("D13").Value = "'0.09149"
$ws.Range("E13").Value = "'5.28%"
$ws.Range("G13").Value = "'21"

$ws.Range("D14").Value = "'0.03345"
$ws.Range("E14").Value = "'-1.49%"
$ws.Range("G14").Value = "'21"

$ws.Range("D15").Value = "'0.09613"
$ws.Range("E15").Value = "'-0.92%"
$ws.Range("G15").Value = "'21"

$ws.Range("D16").Value = "'0.001378"
$ws.Range("E16").Value = "'1.07%"
$ws.Range("G16").Value = "'21"

$ws.Range("D17").Value = "'0.005758"
$ws.Range("E17").Value = "'-1.88%"
$ws.Range("G17").Value = "'21"

$ws.Range("E18").Value = "'-1.91%"
$ws.Range("G18").Value = "'21"

$ws.Range("D19").Value = "'0.3445"
$ws.Range("E19").Value = "'1.14%"
$ws.Range("G19").Value = "'21"

$ws.Range("D20").Value = "'5.261"
$ws.Range("E20").Value = "'4.83%"
$ws.Range("G20").Value = "'21"

$ws.Range("E21").Value = "'-0.21%"
$ws.Range("G21").Value = "'21"

$ws.Range("E22").Value = "'4.08%"
$ws.Range("G22").Value = "'21"

$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04367"
$ws.Range("E23").Value = "'1.08%"
$ws.Range("G23").Value = "'21"

$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001250"
$ws.Range("E24").Value = "'3.02%"
$ws.Range("G24").Value = "'21"

$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "'0.004713"
$ws.Range("E25").Value = "'5.13%"
$ws.Range("G25").Value = "'21"

$ws.Range("B26").Value = "NitroEx"
$ws.Range("C26").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D26").Value = "'0.0001222"
$ws.Range("E26").Value = "'-9.62%"
$ws.Range("G26").Value = "'21"

$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D27").Value = "'0.0003990"
$ws.Range("E27").Value = "'-98.11%"
$ws.Range("G27").Value = "'21"

$ws.Range("G28").Value = "'21"

$ws.Range("G29").Value = "'21"

$ws.Range("G30").Value = "'21"

$ws.Range("G31").Value = "'21"

$ws.Range("G32").Value = "'21"

$ws.Range("G33").Value = "'21"

$ws.Range("G34").Value = "'21"

$ws.Range("G35").Value = "'21"

$ws.Range("G36").Value = "'21"

$ws.Range("G37").Value = "'21"

$ws.Range("G38").Value = "'21"

$ws.Range("D39").Value = "'0.02330"
$ws.Range("E39").Value = "'4.40%"
$ws.Range("G39").Value = "'21"

$ws.Range("D40").Value = "'0.05120"
$ws.Range("E40").Value = "'4.67%"
$ws.Range("G40").Value = "'21"

$ws.Range("D41").Value = "'0.007464"
$ws.Range("E41").Value = "'-1.30%"
$ws.Range("G41").Value = "'21"

$ws.Range("D42").Value = "'0.009032"
$ws.Range("E42").Value = "'-7.70%"
$ws.Range("G42").Value = "'21"

$ws.Range("D43").Value = "'0.1359"
$ws.Range("E43").Value = "'1.82%"
$ws.Range("G43").Value = "'21"

$ws.Range("D44").Value = "'0.002003"
$ws.Range("E44").Value = "'0.31%"
$ws.Range("G44").Value = "'21"

$ws.Range("D45").Value = "'0.008626"
$ws.Range("E45").Value = "'-2.11%"
$ws.Range("G45").Value = "'21"

$ws.Range("D46").Value = "'0.00006629"
$ws.Range("E46").Value = "'0.57%"
$ws.Range("G46").Value = "'21"

$ws.Range("E47").Value = "'-0.17%"
$ws.Range("G47").Value = "'21"

$ws.Range("D48").Value = "'0.003362"
$ws.Range("E48").Value = "'12.06%"
$ws.Range("G48").Value = "'21"

$ws.Range("E49").Value = "'-7.71%"
$ws.Range("G49").Value = "'21"

$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("G50").Value = "'21"

$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("G51").Value = "'21"
